# Update the division-fact worksheet table with new problems.
# The table has 20 rows x 5 columns, but only every 4th row (1, 5, 9, 13, 17)
# holds the visible "a÷b=" expressions; addressing cells by (row, col)
# avoids ambiguity from duplicate text like "70÷8=" appearing twice.

$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "12÷9="
$t.Cell(1, 2).Range.Text = "31÷8="
$t.Cell(1, 3).Range.Text = "65÷5="
$t.Cell(1, 4).Range.Text = "78÷2="
$t.Cell(1, 5).Range.Text = "43÷7="

$t.Cell(5, 1).Range.Text = "80÷8="
$t.Cell(5, 2).Range.Text = "83÷8="
$t.Cell(5, 3).Range.Text = "85÷5="
$t.Cell(5, 4).Range.Text = "58÷9="
$t.Cell(5, 5).Range.Text = "97÷7="

$t.Cell(9, 1).Range.Text = "21÷7="
$t.Cell(9, 2).Range.Text = "21÷7="
$t.Cell(9, 3).Range.Text = "45÷8="
$t.Cell(9, 4).Range.Text = "18÷8="
$t.Cell(9, 5).Range.Text = "61÷3="

$t.Cell(13, 1).Range.Text = "35÷7="
$t.Cell(13, 2).Range.Text = "56÷4="
$t.Cell(13, 3).Range.Text = "92÷2="
$t.Cell(13, 4).Range.Text = "25÷4="
$t.Cell(13, 5).Range.Text = "88÷7="

$t.Cell(17, 1).Range.Text = "47÷5="
$t.Cell(17, 2).Range.Text = "91÷2="
$t.Cell(17, 3).Range.Text = "75÷7="
$t.Cell(17, 4).Range.Text = "21÷4="
$t.Cell(17, 5).Range.Text = "86÷7="
